$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates (Ligand/Receptor average & total expression values, and derived specificities)
$ws.Range("G2").Value = 0.1347866666666667
$ws.Range("H2").Value = 0.40436
$ws.Range("I2").Value = 0.03419045085634245
$ws.Range("J2").Value = 0.03419045085634244
$ws.Range("Q2").Value = 0.03961824929333333
$ws.Range("R2").Value = 0.35656424364
$ws.Range("S2").Value = 0.03419045085634245
$ws.Range("T2").Value = 0.03419045085634244

# Row 3 updates (only derived specificity columns changed)
$ws.Range("I3").Value = 0.3318597741685039
$ws.Range("J3").Value = 0.3318597741685039
$ws.Range("S3").Value = 0.3318597741685039
$ws.Range("T3").Value = 0.3318597741685039

# Row 4 updates (only derived specificity columns changed)
$ws.Range("I4").Value = 0.6339497749751537
$ws.Range("J4").Value = 0.6339497749751537
$ws.Range("S4").Value = 0.6339497749751537
$ws.Range("T4").Value = 0.6339497749751537
